$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.905.44"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.815.38"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'309.00"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4645"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'0.8694"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "'20.29"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.873.93"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("D13").Value = "'5.376"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'0.07093"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "'6.509"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'0.000008711"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'14.65"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "26.952.55"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'10.62"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "2.051.39"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "'1.894"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'150.84"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "'18.32"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").Value = "'5.266"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").Value = "'115.07"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "'0.08897"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'1.158"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "'4.485"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'0.05283"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").Value = "'0.01947"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "'2.986"
$ws.Range("E40").Value = "  +2.09%  "
$ws.Range("D41").Value = "'7.238"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'0.5298"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'2.288"
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'8.433"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("D49").Value = "'103.32"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("D51").Value = "'0.06293"
$ws.Range("E51").Value = "  +0.13%  "
